# Daily auto push: a new timestamp row was logged for 2026/02/02 05:00,
# inserted right after the existing "2026/02/02" block (row 755) and
# before the "2026/12/29" block, pushing every following row down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 756 — rows 756:797 shift down to 757:798.
$ws.Rows.Item(756).Insert()

# Column A holds plain text dates (e.g. "2026/12/29"), not real Excel
# dates. Assigning a date-shaped string directly would get auto-converted
# to a date serial number, so force text formatting first, then drop the
# explicit number-format style again (Style = "Normal") so the cell ends
# up as an unstyled text cell, matching the rest of the column.
$ws.Cells.Item(756, 1).NumberFormat = "@"
$ws.Cells.Item(756, 1).Value = "2026/02/02"
$ws.Cells.Item(756, 1).Style = "Normal"

$ws.Cells.Item(756, 2).Value = "月"
$ws.Cells.Item(756, 3).Value = 5
$ws.Cells.Item(756, 4).Value = 201
